$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = -0.2762176144268791
$ws1.Range("C2").Value = -0.2488306182889404
$ws1.Range("B3").Value = -0.4599053708787398
$ws1.Range("C3").Value = -0.3910610995737829
$ws1.Range("B4").Value = -0.8895767500632773
$ws1.Range("C4").Value = -0.9916701249079565

$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -0.5496112782640099
$ws2.Range("C2").Value = 0.02537908584381716
$ws2.Range("B3").Value = -0.8152621695976044
$ws2.Range("C3").Value = -0.5800336624975728
$ws2.Range("B4").Value = -1.166211979426532
$ws2.Range("C4").Value = 0.4814866994593548
